# Updated symbol list on Fri Jan 20 21:10:10 UTC 2023 with GitHub Actions
# Refresh crypto price/volume/hour snapshot values in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'1.45%"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").Value = "'21"
$ws.Range("G2").ClearFormats()

$ws.Range("D3").Value = "'31.78"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'2.39%"
$ws.Range("E3").ClearFormats()
$ws.Range("G3").Value = "'21"
$ws.Range("G3").ClearFormats()

$ws.Range("D4").Value = "'4.996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'1.37%"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").Value = "'21"
$ws.Range("G4").ClearFormats()

$ws.Range("D5").Value = "'0.07720"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'5.05%"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").Value = "'21"
$ws.Range("G5").ClearFormats()

$ws.Range("D6").Value = "'2.235"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-2.55%"
$ws.Range("E6").ClearFormats()
$ws.Range("G6").Value = "'21"
$ws.Range("G6").ClearFormats()

$ws.Range("D7").Value = "'7.904"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'2.15%"
$ws.Range("E7").ClearFormats()
$ws.Range("G7").Value = "'21"
$ws.Range("G7").ClearFormats()

$ws.Range("D8").Value = "'0.9246"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.75%"
$ws.Range("E8").ClearFormats()
$ws.Range("G8").Value = "'21"
$ws.Range("G8").ClearFormats()

$ws.Range("D9").Value = "'0.09796"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'22.45%"
$ws.Range("E9").ClearFormats()
$ws.Range("G9").Value = "'21"
$ws.Range("G9").ClearFormats()

$ws.Range("D10").Value = "'0.1749"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'3.60%"
$ws.Range("E10").ClearFormats()
$ws.Range("G10").Value = "'21"
$ws.Range("G10").ClearFormats()

$ws.Range("D11").Value = "'0.08432"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'3.13%"
$ws.Range("E11").ClearFormats()
$ws.Range("G11").Value = "'21"
$ws.Range("G11").ClearFormats()

$ws.Range("D12").Value = "'0.03264"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'5.17%"
$ws.Range("E12").ClearFormats()
$ws.Range("G12").Value = "'21"
$ws.Range("G12").ClearFormats()

$ws.Range("D13").Value = "'0.09874"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-2.09%"
$ws.Range("E13").ClearFormats()
$ws.Range("G13").Value = "'21"
$ws.Range("G13").ClearFormats()

$ws.Range("D14").Value = "'0.001473"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-2.27%"
$ws.Range("E14").ClearFormats()
$ws.Range("G14").Value = "'21"
$ws.Range("G14").ClearFormats()

$ws.Range("D15").Value = "'0.005729"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-2.02%"
$ws.Range("E15").ClearFormats()
$ws.Range("G15").Value = "'21"
$ws.Range("G15").ClearFormats()

$ws.Range("D16").Value = "'3.533"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'1.55%"
$ws.Range("E16").ClearFormats()
$ws.Range("G16").Value = "'21"
$ws.Range("G16").ClearFormats()

$ws.Range("D17").Value = "'3.803"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'1.51%"
$ws.Range("E17").ClearFormats()
$ws.Range("G17").Value = "'21"
$ws.Range("G17").ClearFormats()

$ws.Range("D18").Value = "'2.177"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'4.91%"
$ws.Range("E18").ClearFormats()
$ws.Range("G18").Value = "'21"
$ws.Range("G18").ClearFormats()

$ws.Range("D19").Value = "'0.3367"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'1.20%"
$ws.Range("E19").ClearFormats()
$ws.Range("G19").Value = "'21"
$ws.Range("G19").ClearFormats()

$ws.Range("D20").Value = "'0.1324"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'1.56%"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").Value = "'21"
$ws.Range("G20").ClearFormats()

$ws.Range("D21").Value = "'4.061"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'2.06%"
$ws.Range("E21").ClearFormats()
$ws.Range("G21").Value = "'21"
$ws.Range("G21").ClearFormats()

$ws.Range("D22").Value = "'0.2275"
$ws.Range("D22").ClearFormats()
$ws.Range("G22").Value = "'21"
$ws.Range("G22").ClearFormats()

$ws.Range("D23").Value = "'0.04515"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.71%"
$ws.Range("E23").ClearFormats()
$ws.Range("G23").Value = "'21"
$ws.Range("G23").ClearFormats()

$ws.Range("D24").Value = "'0.001213"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'0.19%"
$ws.Range("E24").ClearFormats()
$ws.Range("G24").Value = "'21"
$ws.Range("G24").ClearFormats()

$ws.Range("E25").Value = "'-6.21%"
$ws.Range("E25").ClearFormats()
$ws.Range("G25").Value = "'21"
$ws.Range("G25").ClearFormats()

$ws.Range("D26").Value = "'0.0001288"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'-1.01%"
$ws.Range("E26").ClearFormats()
$ws.Range("G26").Value = "'21"
$ws.Range("G26").ClearFormats()

$ws.Range("D27").Value = "'0.0003364"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-1.01%"
$ws.Range("E27").ClearFormats()
$ws.Range("G27").Value = "'21"
$ws.Range("G27").ClearFormats()

$ws.Range("G28").Value = "'21"
$ws.Range("G28").ClearFormats()

$ws.Range("G29").Value = "'21"
$ws.Range("G29").ClearFormats()

$ws.Range("G30").Value = "'21"
$ws.Range("G30").ClearFormats()

$ws.Range("G31").Value = "'21"
$ws.Range("G31").ClearFormats()

$ws.Range("G32").Value = "'21"
$ws.Range("G32").ClearFormats()

$ws.Range("G33").Value = "'21"
$ws.Range("G33").ClearFormats()

$ws.Range("G34").Value = "'21"
$ws.Range("G34").ClearFormats()

$ws.Range("G35").Value = "'21"
$ws.Range("G35").ClearFormats()

$ws.Range("G36").Value = "'21"
$ws.Range("G36").ClearFormats()

$ws.Range("G37").Value = "'21"
$ws.Range("G37").ClearFormats()

$ws.Range("G38").Value = "'21"
$ws.Range("G38").ClearFormats()

$ws.Range("D39").Value = "'0.01713"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'6.85%"
$ws.Range("E39").ClearFormats()
$ws.Range("G39").Value = "'21"
$ws.Range("G39").ClearFormats()

$ws.Range("D40").Value = "'0.04651"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'4.62%"
$ws.Range("E40").ClearFormats()
$ws.Range("G40").Value = "'21"
$ws.Range("G40").ClearFormats()

$ws.Range("D41").Value = "'0.007704"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'4.62%"
$ws.Range("E41").ClearFormats()
$ws.Range("G41").Value = "'21"
$ws.Range("G41").ClearFormats()

$ws.Range("D42").Value = "'0.009743"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'12.81%"
$ws.Range("E42").ClearFormats()
$ws.Range("G42").Value = "'21"
$ws.Range("G42").ClearFormats()

$ws.Range("D43").Value = "'0.1393"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'4.45%"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").Value = "'21"
$ws.Range("G43").ClearFormats()

$ws.Range("D44").Value = "'0.001982"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'2.05%"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").Value = "'21"
$ws.Range("G44").ClearFormats()

$ws.Range("D45").Value = "'0.009714"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'2.04%"
$ws.Range("E45").ClearFormats()
$ws.Range("G45").Value = "'21"
$ws.Range("G45").ClearFormats()

$ws.Range("D46").Value = "'0.00006053"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'1.61%"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").Value = "'21"
$ws.Range("G46").ClearFormats()

$ws.Range("D47").Value = "'0.00000000743"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-1.01%"
$ws.Range("E47").ClearFormats()
$ws.Range("G47").Value = "'21"
$ws.Range("G47").ClearFormats()

$ws.Range("G48").Value = "'21"
$ws.Range("G48").ClearFormats()

$ws.Range("D49").Value = "'0.001982"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'-31.63%"
$ws.Range("E49").ClearFormats()
$ws.Range("G49").Value = "'21"
$ws.Range("G49").ClearFormats()

$ws.Range("D50").Value = "'0.00002081"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-1.01%"
$ws.Range("E50").ClearFormats()
$ws.Range("G50").Value = "'21"
$ws.Range("G50").ClearFormats()

$ws.Range("D51").Value = "'0.0001982"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'-1.01%"
$ws.Range("E51").ClearFormats()
$ws.Range("G51").Value = "'21"
$ws.Range("G51").ClearFormats()
